# Update diagrams in Developer Guide
#
# The "Spare Parser" box in the Logic Component Class Diagram is renamed
# to "WishBook Parser" - i.e. the first paragraph of shape id=16
# ("Rectangle 62") changes its text from "Spare" to "WishBook". The
# second paragraph ("Parser") is left untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape id=16 ("Rectangle 62") holds the two-line label "Spare" / "Parser".
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Id -eq 16) {
        $shp = $candidate
        break
    }
}

$tr = $shp.TextFrame.TextRange
$firstPara = $tr.Paragraphs(1, 1)
$firstPara.Text = "WishBook"
